$wb = $excel.ActiveWorkbook

# --- "About" sheet: add a new "Notes" section (rows 48-50) explaining the
#     rebound-effect meaning of the elasticity values used elsewhere. ---
$about = $wb.Worksheets.Item("About")

[void]($about.Range("A48").Font.Bold = $true)
$about.Range("A48").Value = "Notes"
$about.Range("A49").Value = "This variable is also known as the ""Fuel Economy Rebound Effect"" or ""Fuel Cost Rebound Effect."" It is the change"
$about.Range("A50").Value = "in VMT as a fraction of the change in fuel cost. E.g. for a 1% increase in fuel cost per mile, VMT changes by -0.1%."

# --- "EoDfVUwFC" sheet: clarify header label, wrap it, and give the row more
#     height so the wrapped text is fully visible. ---
$eod = $wb.Worksheets.Item("EoDfVUwFC")

[void]($eod.Range("B1").WrapText = $true)
$eod.Range("B1").Value = "Elasticity (dimensionless)"
$eod.Rows.Item(1).RowHeight = 30

# Leave the same on-screen selection/scroll state the author ended up with:
# EoDfVUwFC's header cell selected, but "About" remains the active tab.
[void]$eod.Activate()
[void]$eod.Range("B1").Select()

[void]$about.Activate()
[void]$about.Range("A48").Select()
